$wb = $excel.ActiveWorkbook

# ---------- Sheet 1: Summary ----------
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B2").Value = 0.602996254681648
$wsSummary.Range("C2").Value = 0.5625
$wsSummary.Range("D2").Value = 0.9269662921348315
$wsSummary.Range("E2").Value = 0.7001414427157001
$wsSummary.Range("F2").Value = 0.8206233421750663
$wsSummary.Range("G2").Value = 0.9044272663387211
$wsSummary.Range("H2").Value = 0.7682777146544348
$wsSummary.Range("I2").Value = 495
$wsSummary.Range("J2").Value = 385
$wsSummary.Range("K2").Value = 149
$wsSummary.Range("L2").Value = 39

# ---------- Sheet 2: Classification Report ----------
$wsClassRep = $wb.Worksheets.Item("Classification Report")

# Row 2 - class "0"
$wsClassRep.Range("B2").Value = 0.7925531914893617
$wsClassRep.Range("C2").Value = 0.2790262172284644
$wsClassRep.Range("D2").Value = 0.4127423822714681

# Row 3 - class "1"
$wsClassRep.Range("B3").Value = 0.5625
$wsClassRep.Range("C3").Value = 0.9269662921348315
$wsClassRep.Range("D3").Value = 0.7001414427157001

# Row 4 - accuracy
$wsClassRep.Range("B4").Value = 0.602996254681648
$wsClassRep.Range("C4").Value = 0.602996254681648
$wsClassRep.Range("D4").Value = 0.602996254681648
$wsClassRep.Range("E4").Value = 0.602996254681648

# Row 5 - macro avg
$wsClassRep.Range("B5").Value = 0.6775265957446808
$wsClassRep.Range("C5").Value = 0.602996254681648
$wsClassRep.Range("D5").Value = 0.5564419124935841

# Row 6 - weighted avg
$wsClassRep.Range("B6").Value = 0.6775265957446808
$wsClassRep.Range("C6").Value = 0.602996254681648
$wsClassRep.Range("D6").Value = 0.5564419124935841

# ---------- Sheet 3: Confusion Matrix ----------
$wsConfMat = $wb.Worksheets.Item("Confusion Matrix")

# Row 2 - Actual 0
$wsConfMat.Range("B2").Value = 149
$wsConfMat.Range("C2").Value = 385

# Row 3 - Actual 1
$wsConfMat.Range("B3").Value = 39
$wsConfMat.Range("C3").Value = 495
